# Apply diary bot updates:
#  - Diary sheet: drop the stray empty F10 cell, append 3 new diary_entry rows (11-13)
#  - Safety sheet: append 1 new crisis row (4) for the suicidal-thoughts entry logged in row 12

$wb = $excel.ActiveWorkbook
$diary = $wb.Worksheets.Item("Diary")
$safety = $wb.Worksheets.Item("Safety")

# 1) F10 was an empty placeholder cell in the source data; remove it entirely.
$diary.Range("F10").ClearContents()

# 2) Row 11 - new diary entry, no crisis detected.
$diary.Cells.Item(11, 1).Value = 539011121
$diary.Cells.Item(11, 2).Value = "Anna_Safonova_life"
$diary.Cells.Item(11, 3).Value = "Анна"
$diary.Cells.Item(11, 4).Value = "diary_entry"
$diary.Cells.Item(11, 5).Value = "Чувствую себя хорошо, но ещё есть тревога"
$diary.Cells.Item(11, 7).Value = "2025-10-19 12:34:03"

# 3) Row 12 - new diary entry, suicidal-thoughts crisis detected.
$diary.Cells.Item(12, 1).Value = 6479033897
$diary.Cells.Item(12, 2).Value = "zhuravlstrogo"
$diary.Cells.Item(12, 3).Value = "User"
$diary.Cells.Item(12, 4).Value = "diary_entry"
$diary.Cells.Item(12, 5).Value = "Хочу умереть"
$diary.Cells.Item(12, 7).Value = "2025-10-19 12:51:36"

# 4) Row 13 - new diary entry, no crisis detected. The upstream bot left its
#    usual stray empty Date-Time cell behind (same glitch fixed in F10 above).
$diary.Cells.Item(13, 1).Value = 6893133357
$diary.Cells.Item(13, 2).Value = "nadzh_k"
$diary.Cells.Item(13, 3).Value = "User"
$diary.Cells.Item(13, 4).Value = "diary_entry"
$diary.Cells.Item(13, 5).Value = "Хочу начать любить свою жизнь"
$diary.Cells.Item(13, 6).WrapText = $false
$diary.Cells.Item(13, 7).Value = "2025-10-19 15:19:09"

# 5) Log the crisis found in row 12 to the Safety sheet.
$safety.Cells.Item(4, 1).Value = 6479033897
$safety.Cells.Item(4, 2).Value = "zhuravlstrogo"
$safety.Cells.Item(4, 3).Value = "2025-10-19 12:51:36"
$safety.Cells.Item(4, 4).Value = "Суицидальные мысли"
$safety.Cells.Item(4, 5).Value = "diary"
$safety.Cells.Item(4, 6).Value = "Хочу умереть"

$wb.Save()
